# This script applies a weekly data refresh to the "Mandarina" price sheet.
# It inserts two new rows (a new week of observations for variety "Murcott")
# right above the existing block of price rows, shifting all subsequent rows
# down by two (dimension grows from A1:T153 to A1:T155).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 62; this shifts the existing rows 62-153
# down to 64-155 and keeps the column formatting (including the date format
# on column D) that Excel propagates from the row immediately below.
$ws.Rows("62:63").Insert()

# --- New row 62: Mandarina / Murcott / Primera -----------------------------
$ws.Range("A62").Value = 7
$ws.Range("B62").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C62").Value = "Ñuble"
$ws.Range("D62").Value = 44483
$ws.Range("E62").Value = 16
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100102
$ws.Range("H62").Value = "Cítricos"
$ws.Range("I62").Value = 100102004
$ws.Range("J62").Value = "Mandarina"
$ws.Range("K62").Value = "Murcott"
$ws.Range("L62").Value = "Primera"
$ws.Range("M62").Value = 240
$ws.Range("N62").Value = 5500
$ws.Range("O62").Value = 6000
$ws.Range("P62").Value = 5750
$ws.Range("Q62").Value = "$/bandeja 10 kilos"
$ws.Range("R62").Value = "Provincia de Limarí"
$ws.Range("S62").Value = 575
$ws.Range("T62").Value = 10

# --- New row 63: Mandarina / Murcott / Segunda ------------------------------
$ws.Range("A63").Value = 7
$ws.Range("B63").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C63").Value = "Ñuble"
$ws.Range("D63").Value = 44483
$ws.Range("E63").Value = 16
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100102
$ws.Range("H63").Value = "Cítricos"
$ws.Range("I63").Value = 100102004
$ws.Range("J63").Value = "Mandarina"
$ws.Range("K63").Value = "Murcott"
$ws.Range("L63").Value = "Segunda"
$ws.Range("M63").Value = 240
$ws.Range("N63").Value = 4500
$ws.Range("O63").Value = 5000
$ws.Range("P63").Value = 4750
$ws.Range("Q63").Value = "$/bandeja 10 kilos"
$ws.Range("R63").Value = "Provincia de Limarí"
$ws.Range("S63").Value = 475
$ws.Range("T63").Value = 10
